$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 176.125
$ws.Range("I33").Value = 176.125
$ws.Range("K33").Value = 176.125
$ws.Range("M33").Value = 52.875
$ws.Range("H57").Value = 63571.125
$ws.Range("J57").Value = 63571.125
$ws.Range("L57").Value = 190713.375
$ws.Range("N57").Value = -191711.375
$ws.Range("H95").Value = 34749.332
$ws.Range("J95").Value = 34749.332
$ws.Range("L95").Value = 34749.332
$ws.Range("N95").Value = -40241.332
$ws.Range("H103").Value = 83334590
$ws.Range("I103").Value = 500
$ws.Range("J103").Value = 100001410
$ws.Range("K103").Value = 1500
$ws.Range("L103").Value = 300004230
$ws.Range("M103").Value = -914
$ws.Range("N103").Value = -300005402
$ws.Range("H132").Value = 1355.1719
$ws.Range("I132").Value = 1375.661
$ws.Range("J132").Value = 1113.4
$ws.Range("K132").Value = 4126.983
$ws.Range("L132").Value = 3340.2
$ws.Range("M132").Value = -1596.983
$ws.Range("N132").Value = -8400.200000000001
$ws.Range("H137").Value = 1242.909
$ws.Range("I137").Value = 1151.625
$ws.Range("J137").Value = 1486.3334
$ws.Range("K137").Value = 3454.875
$ws.Range("L137").Value = 4459.0002
$ws.Range("M137").Value = -904.875
$ws.Range("N137").Value = -9559.0002
$ws.Range("H141").Value = 2895.6365
$ws.Range("I141").Value = 3242.75
$ws.Range("J141").Value = 1970
$ws.Range("K141").Value = 9728.25
$ws.Range("L141").Value = 5910
$ws.Range("M141").Value = -4548.25
$ws.Range("N141").Value = -16270

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4996.8223
$ws.Range("I32").Value = 4424.325
$ws.Range("J32").Value = 9576.799999999999
$ws.Range("K32").Value = 4424.325
$ws.Range("L32").Value = 9576.799999999999
$ws.Range("M32").Value = -4137.325
$ws.Range("N32").Value = -10150.8
$ws.Range("H57").Value = 5155
$ws.Range("I57").Value = 5155
$ws.Range("K57").Value = 5155
$ws.Range("M57").Value = -4671
$ws.Range("H61").Value = 1770.2
$ws.Range("I61").Value = 1287.0714
$ws.Range("K61").Value = 1287.0714
$ws.Range("M61").Value = -1075.0714
$ws.Range("H74").Value = 2149.6667
$ws.Range("I74").Value = 2149.6667
$ws.Range("K74").Value = 2149.6667
$ws.Range("M74").Value = -1275.6667
$ws.Range("H77").Value = 2149.6667
$ws.Range("I77").Value = 2149.6667
$ws.Range("K77").Value = 10748.3335
$ws.Range("M77").Value = -6380.333500000001
$ws.Range("H95").Value = 30204
$ws.Range("J95").Value = 30204
$ws.Range("L95").Value = 30204
$ws.Range("N95").Value = -35696
$ws.Range("H97").Value = 1560.5454
$ws.Range("I97").Value = 1560.5454
$ws.Range("K97").Value = 1560.5454
$ws.Range("M97").Value = -1064.5454
$ws.Range("H102").Value = 4060.5557
$ws.Range("I102").Value = 4060.5557
$ws.Range("K102").Value = 4060.5557
$ws.Range("M102").Value = -2438.5557
$ws.Range("H133").Value = 98761
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("H136").Value = 1770.2
$ws.Range("I136").Value = 1287.0714
$ws.Range("K136").Value = 3861.2142
$ws.Range("M136").Value = -1311.2142

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 897
$ws.Range("I86").Value = 885.2308
$ws.Range("J86").Value = 922.5
$ws.Range("K86").Value = 885.2308
$ws.Range("L86").Value = 922.5
$ws.Range("M86").Value = 237.7692
$ws.Range("N86").Value = -3168.5
$ws.Range("H89").Value = 897
$ws.Range("I89").Value = 885.2308
$ws.Range("J89").Value = 922.5
$ws.Range("K89").Value = 4426.154
$ws.Range("L89").Value = 4612.5
$ws.Range("M89").Value = 1189.846
$ws.Range("N89").Value = -15844.5
$ws.Range("H94").Value = 1814.2307
$ws.Range("I94").Value = 1814.2307
$ws.Range("K94").Value = 1814.2307
$ws.Range("M94").Value = -1363.2307
$ws.Range("H105").Value = 3938
$ws.Range("I105").Value = 7865.4287
$ws.Range("J105").Value = 883.3333
$ws.Range("K105").Value = 7865.4287
$ws.Range("L105").Value = 883.3333
$ws.Range("M105").Value = -6118.4287
$ws.Range("N105").Value = -4377.3333
$ws.Range("H113").Value = 5113
$ws.Range("I113").Value = 5113
$ws.Range("K113").Value = 5113
$ws.Range("M113").Value = -2943
$ws.Range("H128").Value = 4395.857
$ws.Range("I128").Value = 4395.857
$ws.Range("K128").Value = 13187.571
$ws.Range("M128").Value = -10697.571
$ws.Range("H134").Value = 1687.6666
$ws.Range("I134").Value = 1687.6666
$ws.Range("K134").Value = 5062.9998
$ws.Range("M134").Value = -2527.9998

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H7").Value = 45457404
$ws.Range("I7").Value = 66668764
$ws.Range("J7").Value = 4492.857
$ws.Range("K7").Value = 66668764
$ws.Range("L7").Value = 4492.857
$ws.Range("M7").Value = -66668651
$ws.Range("N7").Value = -4718.857
$ws.Range("H31").Value = 15092.565
$ws.Range("I31").Value = 1614.7778
$ws.Range("J31").Value = 63612.6
$ws.Range("K31").Value = 1614.7778
$ws.Range("L31").Value = 63612.6
$ws.Range("M31").Value = -1319.7778
$ws.Range("N31").Value = -64202.6
$ws.Range("H34").Value = 15092.565
$ws.Range("I34").Value = 1614.7778
$ws.Range("J34").Value = 63612.6
$ws.Range("K34").Value = 1614.7778
$ws.Range("L34").Value = 63612.6
$ws.Range("M34").Value = -1412.7778
$ws.Range("N34").Value = -64016.6
$ws.Range("H41").Value = 10150
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H50").Value = 15000
$ws.Range("J50").Value = 15000
$ws.Range("L50").Value = 15000
$ws.Range("N50").Value = -16250
$ws.Range("H51").Value = 15000
$ws.Range("J51").Value = 15000
$ws.Range("L51").Value = 15000
$ws.Range("N51").Value = -16472
$ws.Range("H59").Value = 19998.334
$ws.Range("H61").Value = 15000
$ws.Range("J61").Value = 15000
$ws.Range("L61").Value = 15000
$ws.Range("N61").Value = -15696
$ws.Range("H76").Value = 5073.3335
$ws.Range("I76").Value = 5073.3335
$ws.Range("K76").Value = 5073.3335
$ws.Range("M76").Value = -4758.3335
$ws.Range("H79").Value = 5073.3335
$ws.Range("I79").Value = 5073.3335
$ws.Range("K79").Value = 5073.3335
$ws.Range("M79").Value = -3981.3335
$ws.Range("H105").Value = 1616.1666
$ws.Range("J105").Value = 1451.8
$ws.Range("L105").Value = 1451.8
$ws.Range("N105").Value = -4945.8
$ws.Range("H134").Value = 1830.6129
$ws.Range("I134").Value = 1267.2069
$ws.Range("K134").Value = 3801.620699999999
$ws.Range("M134").Value = -1266.620699999999

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 114.42857
$ws.Range("I7").Value = 141.8
$ws.Range("J7").Value = 46
$ws.Range("K7").Value = 425.4
$ws.Range("L7").Value = 138
$ws.Range("M7").Value = -313.4
$ws.Range("N7").Value = -362
$ws.Range("H15").Value = 65.27273
$ws.Range("I15").Value = 48.666668
$ws.Range("J15").Value = 140
$ws.Range("K15").Value = 146.000004
$ws.Range("L15").Value = 420
$ws.Range("M15").Value = -6.00000399999999
$ws.Range("N15").Value = -700
$ws.Range("H80").Value = 3540.3
$ws.Range("I80").Value = 1500
$ws.Range("J80").Value = 3767
$ws.Range("K80").Value = 4500
$ws.Range("L80").Value = 11301
$ws.Range("M80").Value = -3564
$ws.Range("N80").Value = -13173
$ws.Range("H83").Value = 3540.3
$ws.Range("I83").Value = 1500
$ws.Range("J83").Value = 3767
$ws.Range("K83").Value = 13500
$ws.Range("L83").Value = 33903
$ws.Range("M83").Value = -8820
$ws.Range("N83").Value = -43263
$ws.Range("H117").Value = 7090.077
$ws.Range("J117").Value = 7347.5835
$ws.Range("L117").Value = 22042.7505
$ws.Range("N117").Value = -28926.7505
$ws.Range("H131").Value = 1681.5
$ws.Range("I131").Value = 1312.8889
$ws.Range("J131").Value = 4999
$ws.Range("K131").Value = 3938.6667
$ws.Range("L131").Value = 14997
$ws.Range("M131").Value = 1101.3333
$ws.Range("N131").Value = -25077

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2140.963
$ws.Range("I132").Value = 2194.4614
$ws.Range("J132").Value = 750
$ws.Range("K132").Value = 6583.3842
$ws.Range("L132").Value = 2250
$ws.Range("M132").Value = -4053.3842
$ws.Range("N132").Value = -7310

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H60").Value = 59999
$ws.Range("J60").Value = 59999
$ws.Range("L60").Value = 59999
$ws.Range("N60").Value = -61017
$ws.Range("H82").Value = 1897.4375
$ws.Range("J82").Value = 1991.8572
$ws.Range("L82").Value = 1991.8572
$ws.Range("N82").Value = -2713.8572
$ws.Range("H85").Value = 1897.4375
$ws.Range("J85").Value = 1991.8572
$ws.Range("L85").Value = 1991.8572
$ws.Range("N85").Value = -4487.8572
$ws.Range("H136").Value = 5022.1
$ws.Range("I136").Value = 3204.3333
$ws.Range("J136").Value = 7748.75
$ws.Range("K136").Value = 9612.999899999999
$ws.Range("L136").Value = 23246.25
$ws.Range("M136").Value = -7062.999899999999
$ws.Range("N136").Value = -28346.25

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 47516.75
$ws.Range("I54").Value = 40070
$ws.Range("K54").Value = 40070
$ws.Range("M54").Value = -39550
